$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Edit 1: Title-page author line "Antti Minkkinen G7786" -> split
# into separate runs with spellcheck proofErr markers around the
# two Finnish name tokens.
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Antti Minkkinen G7786")
$para1 = $rng.Paragraphs(1).Range
$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>Antti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>Minkkinen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> G7786</w:t></w:r></w:p>
'@
$para1.InsertXML($xml1)

# ---------------------------------------------------------------
# Edit 2: Introduction paragraph - "...is to create a" + "n" ->
# "...is to create " + "a prototype of " + "a" + "n"
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("The idea of the project")
$para2 = $rng.Paragraphs(1).Range
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The idea of the project </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>is</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> to create </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">a prototype of </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>a</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>n</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> application to </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">share </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>the prices of certain products</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> in</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> certain</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> shops</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>and help find the cheapest</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> product</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>from the nearest shops.</w:t></w:r></w:p>
'@
$para2.InsertXML($xml2)

# ---------------------------------------------------------------
# Edit 3: move a lastRenderedPageBreak from the "These activities
# are used to add a price..." paragraph to the picture paragraph
# right before it (relativeHeight 251660288).
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Product name and category (and sub categories) can be selected.")
$picPara = $rng.Paragraphs(1).Next()
$xml3 = @'
<w:p w:rsidR="00830D96" w:rsidRDefault="00E947EE" w:rsidP="00A21C8A"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="fi-FI"/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251660288" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="50F108F8" wp14:editId="49F1E749"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="margin"><wp:align>left</wp:align></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>384810</wp:posOffset></wp:positionV><wp:extent cx="1814249" cy="3457575"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapTopAndBottom/><wp:docPr id="7" name="Picture 7"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="7" name="5 Add price 1.png"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId9" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1814249" cy="3457575"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="page"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="page"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="fi-FI"/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251662336" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="1D311CB6" wp14:editId="00273914"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="margin"><wp:align>right</wp:align></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>323850</wp:posOffset></wp:positionV><wp:extent cx="1828800" cy="3486150"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapTopAndBottom/><wp:docPr id="9" name="Picture 9"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="9" name="5 Add price 3.png"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId10" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1828800" cy="3486150"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="page"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="page"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="fi-FI"/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251661312" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="7BC56843" wp14:editId="5F3DD4CA"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="margin"><wp:align>center</wp:align></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>375285</wp:posOffset></wp:positionV><wp:extent cx="1838325" cy="3503295"/><wp:effectExtent l="0" t="0" r="9525" b="1905"/><wp:wrapTopAndBottom/><wp:docPr id="8" name="Picture 8"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="8" name="5 Add price 2.png"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId11" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1838325" cy="3503295"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="page"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="page"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r></w:p>
'@
$picPara.Range.InsertXML($xml3)

$rng = $d.Content
$rng.Find.Execute("These activities are used to add a price to a product in a certain shop.")
$pricePara = $rng.Paragraphs(1).Range
$xml3b = @'
<w:p w:rsidR="00A95AD9" w:rsidRDefault="00830D96" w:rsidP="00A21C8A"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>These activities are used to add a price to a product in a certain shop.</w:t></w:r><w:r w:rsidR="00E947EE"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> The product can be manually selected or it can be automatically selected if the user navigated here from other product. In the second screen the user selects the shop from the map or by searching from the list. The list automatically contains the nearest shops and also shows the distance to them. The third screen is used for adding the price(s) of the product.</w:t></w:r></w:p>
'@
$pricePara.InsertXML($xml3b)

# ---------------------------------------------------------------
# Edit 4 + 5: picture paragraph (relativeHeight 251663360) gets a
# lastRenderedPageBreak; the trailing sentence is re-split into
# three runs (capitalising GPS); the _GoBack bookmark is removed
# from here (it moves to the new heading added at the end).
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("These activities are used to add a price to a product in a certain shop.")
$shopPara = $rng.Paragraphs(1).Next().Next().Next()
$xml4 = @'
<w:p w:rsidR="004C38E6" w:rsidRPr="00830D96" w:rsidRDefault="0074395B" w:rsidP="004C38E6"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="fi-FI"/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251663360" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>0</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>0</wp:posOffset></wp:positionV><wp:extent cx="1733550" cy="3303781"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapTopAndBottom/><wp:docPr id="6" name="Picture 6"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="6" name="4 Add Shop.png"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId12" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1733550" cy="3303781"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="page"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="page"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r><w:r w:rsidR="00830D96"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>This activity is used to add a shop to the database.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> The user can type the name of the shop and address or they </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>can get the address from map/GPS</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
$shopPara.Range.InsertXML($xml4)

# ---------------------------------------------------------------
# Edit 6: the final empty paragraph becomes a "Work load and
# planning" Heading1 (keeping the _GoBack bookmark) followed by
# four new paragraphs describing the weekly plan.
# ---------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$xml5 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Work load and planning</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Week 45</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Database setup </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">(maybe </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>digital ocean</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>?) and starting to implement business logic and activities.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Week 46</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Database working and some data can be inserted from app.</w:t></w:r></w:p>
'@
$lastPara.Range.InsertXML($xml5)

Write-Output "all edits applied"
